$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("AG2").Value = 6
$ws.Range("AO2").Value = 19
$ws.Range("AR2").Value = 126
$ws.Range("G2").Value = 2.8
$ws.Range("H2").Value = 2.8
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 1.8
$ws.Range("N2").Value = 4.75
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 11
$ws.Range("Z2").Value = 29

# Row 3
$ws.Range("AA3").Value = 21
$ws.Range("AC3").Value = 7
$ws.Range("AF3").Value = 67
$ws.Range("AO3").Value = 13
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 9
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73

# Row 4
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.62

# Row 7
$ws.Range("I7").Value = 2.38

# Row 9
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3

# Row 10
$ws.Range("AA10").Value = 17
$ws.Range("AG10").Value = 10
$ws.Range("AH10").Value = 17
$ws.Range("AI10").Value = 13
$ws.Range("AJ10").Value = 41
$ws.Range("AK10").Value = 29
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 12
$ws.Range("AR10").Value = 51
$ws.Range("AS10").Value = 151
$ws.Range("AW10").Value = 5.5
$ws.Range("AX10").Value = 19
$ws.Range("AY10").Value = 29
$ws.Range("AZ10").Value = 67
$ws.Range("G10").Value = 2.15
$ws.Range("H10").Value = 3.25
$ws.Range("I10").Value = 3.4
$ws.Range("J10").Value = 2.88
$ws.Range("K10").Value = 2.1
$ws.Range("L10").Value = 4
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.4
$ws.Range("Q10").Value = 2.05
$ws.Range("R10").Value = 1.8
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 1.91
$ws.Range("W10").Value = 7.5
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 9
$ws.Range("Z10").Value = 19

# Row 11
$ws.Range("AD11").Value = 7
$ws.Range("AF11").Value = 126
$ws.Range("AG11").Value = 9.5
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 67
$ws.Range("AK11").Value = 51
$ws.Range("AN11").Value = 3.4
$ws.Range("AO11").Value = 10
$ws.Range("AP11").Value = 29
$ws.Range("AQ11").Value = 41
$ws.Range("AU11").Value = 11
$ws.Range("AW11").Value = 7
$ws.Range("G11").Value = 1.73
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 5.75
$ws.Range("J11").Value = 2.5
$ws.Range("K11").Value = 1.83
$ws.Range("L11").Value = 7
$ws.Range("X11").Value = 6.5
$ws.Range("Z11").Value = 13
